# Append the new daily-push row (row 58) to the sheet, extending the
# table from A1:D57 to A1:D58.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 58

# Column A holds a date formatted as plain text (e.g. "2025/10/04"),
# matching the existing rows which store dates as text rather than as
# real Excel date serials. Force the cell to Text format before writing
# the value so Excel doesn't auto-convert the "yyyy/mm/dd"-looking
# string into a date, then clear the format back off so the cell has no
# leftover explicit style (consistent with the other rows in the sheet).
$dateCell = $ws.Range("A" + $newRow)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025/10/04"
$dateCell.ClearFormats()

$ws.Range("B" + $newRow).Value = "土"
$ws.Range("C" + $newRow).Value = 4
$ws.Range("D" + $newRow).Value = 201
